# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rewrites the "Periodo Mora" ledger rows (B16:G29) so that, per worker,
# the period list runs 2108 -> 2102 (descending) instead of 2102 -> 2108
# (ascending), and re-points the "Valor Mora" figures that travel with
# the oldest period (2108) for each worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- ANDRES DE JESUS MARQUEZ LOBO (CC 1047420533) block: rows 16-22 ---
$ws.Range("E16").Value = "2108"
$ws.Range("F16").Value = 33945

$ws.Range("C17").Value = "1047420533"
$ws.Range("D17").Value = "ANDRES DE JESUS MARQUEZ LOBO"
$ws.Range("E17").Value = "2107"

$ws.Range("E18").Value = "2106"

$ws.Range("C19").Value = "1047420533"
$ws.Range("D19").Value = "ANDRES DE JESUS MARQUEZ LOBO"
$ws.Range("E19").Value = "2105"

$ws.Range("C20").Value = "1047420533"
$ws.Range("D20").Value = "ANDRES DE JESUS MARQUEZ LOBO"

$ws.Range("E21").Value = "2103"

$ws.Range("E22").Value = "2102"

# --- JOSE SIMON RHENALS CASSIANI (CC 9284806) block: rows 23-29 ---
$ws.Range("E23").Value = "2108"
$ws.Range("F23").Value = 35129

$ws.Range("C24").Value = "9284806"
$ws.Range("D24").Value = "JOSE SIMON RHENALS CASSIANI"
$ws.Range("E24").Value = "2107"

$ws.Range("C26").Value = "9284806"
$ws.Range("D26").Value = "JOSE SIMON RHENALS CASSIANI"
$ws.Range("E26").Value = "2105"

$ws.Range("E27").Value = "2104"

$ws.Range("E28").Value = "2103"
$ws.Range("F28").Value = 36341

$ws.Range("C29").Value = "9284806"
$ws.Range("D29").Value = "JOSE SIMON RHENALS CASSIANI"
$ws.Range("E29").Value = "2102"
$ws.Range("F29").Value = 36341
